$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-10 from 2023-11-13 (45243) to 2023-11-14 (45244)
$ws.Range("C2:C10").Value = 45244
